$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their original Text storage so that
# numeric-looking strings (e.g. "19.60", "0.0982") are not auto-converted
# to numbers by Excel and lose formatting / change type.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "26.650.65"
$ws.Cells.Item(3, 4).Value = "1.592.48"
$ws.Cells.Item(3, 5).Value = "  -1.79%  "
$ws.Cells.Item(4, 5).Value = "  +0.03%  "
$ws.Cells.Item(5, 4).Value = "211.12"
$ws.Cells.Item(5, 5).Value = "  -1.36%  "
$ws.Cells.Item(6, 5).Value = "  -0.02%  "
$ws.Cells.Item(7, 5).Value = "  +0.01%  "
$ws.Cells.Item(8, 5).Value = "  -1.47%  "
$ws.Cells.Item(9, 5).Value = "  -1.82%  "
$ws.Cells.Item(10, 4).Value = "19.60"
$ws.Cells.Item(10, 5).Value = "  -1.57%  "
$ws.Cells.Item(11, 4).Value = "0.0835"
$ws.Cells.Item(11, 5).Value = "  -0.60%  "
$ws.Cells.Item(12, 4).Value = "1.818.53"
$ws.Cells.Item(12, 5).Value = "  -1.64%  "
$ws.Cells.Item(13, 4).Value = "1.596.18"
$ws.Cells.Item(13, 5).Value = "  -1.80%  "
$ws.Cells.Item(14, 5).Value = "  -2.40%  "
$ws.Cells.Item(15, 5).Value = "  -3.07%  "
$ws.Cells.Item(16, 4).Value = "64.75"
$ws.Cells.Item(16, 5).Value = "  +0.46%  "
$ws.Cells.Item(17, 4).Value = "26.615.53"
$ws.Cells.Item(17, 5).Value = "  -1.58%  "
$ws.Cells.Item(18, 4).Value = "0.0₃0729"
$ws.Cells.Item(18, 5).Value = "  -1.08%  "
$ws.Cells.Item(19, 4).Value = "209.12"
$ws.Cells.Item(19, 5).Value = "  -2.53%  "
$ws.Cells.Item(20, 5).Value = "  +0.05%  "
$ws.Cells.Item(22, 4).Value = "4.23"
$ws.Cells.Item(22, 5).Value = "  -2.31%  "
$ws.Cells.Item(23, 4).Value = "2.29"
$ws.Cells.Item(23, 5).Value = "  -2.01%  "
$ws.Cells.Item(24, 4).Value = "8.89"
$ws.Cells.Item(24, 5).Value = "  -1.48%  "
$ws.Cells.Item(25, 4).Value = "146.50"
$ws.Cells.Item(25, 5).Value = "  -0.63%  "
$ws.Cells.Item(26, 5).Value = "  +0.12%  "
$ws.Cells.Item(27, 5).Value = "  -4.18%  "
$ws.Cells.Item(28, 5).Value = "  -0.17%  "
$ws.Cells.Item(29, 5).Value = "  -1.32%  "
$ws.Cells.Item(30, 5).Value = "  -1.82%  "
$ws.Cells.Item(31, 4).Value = "1.16"
$ws.Cells.Item(31, 5).Value = "  -0.87%  "
$ws.Cells.Item(32, 5).Value = "  -2.85%  "
$ws.Cells.Item(33, 4).Value = "0.687"
$ws.Cells.Item(33, 5).Value = "  -4.50%  "
$ws.Cells.Item(34, 5).Value = "  -3.16%  "
$ws.Cells.Item(35, 4).Value = "1.293.22"
$ws.Cells.Item(35, 5).Value = "  -3.25%  "
$ws.Cells.Item(36, 5).Value = "  -0.54%  "
$ws.Cells.Item(37, 4).Value = "1.47"
$ws.Cells.Item(37, 5).Value = "  -5.22%  "
$ws.Cells.Item(38, 5).Value = "  -2.82%  "
$ws.Cells.Item(39, 5).Value = "  -0.46%  "
$ws.Cells.Item(40, 5).Value = "  +0.10%  "
$ws.Cells.Item(41, 4).Value = "0.789"
$ws.Cells.Item(41, 5).Value = "  -0.66%  "
$ws.Cells.Item(42, 5).Value = "  -1.88%  "
$ws.Cells.Item(43, 4).Value = "5.34"
$ws.Cells.Item(43, 5).Value = "  -0.45%  "
$ws.Cells.Item(44, 4).Value = "63.49"
$ws.Cells.Item(44, 5).Value = "  -0.66%  "
$ws.Cells.Item(45, 4).Value = "1.728.91"
$ws.Cells.Item(45, 5).Value = "  -1.79%  "
$ws.Cells.Item(46, 4).Value = "0.891"
$ws.Cells.Item(46, 5).Value = "  +4.26%  "
$ws.Cells.Item(47, 4).Value = "89.67"
$ws.Cells.Item(47, 5).Value = "  -0.34%  "
$ws.Cells.Item(48, 5).Value = "  -0.49%  "
$ws.Cells.Item(49, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(49, 4).Value = "0.0₆0105"
$ws.Cells.Item(49, 5).Value = "  -1.49%  "
$ws.Cells.Item(50, 2).Value = "Algorand"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(50, 4).Value = "0.0982"
$ws.Cells.Item(50, 5).Value = "  -1.55%  "
$ws.Cells.Item(51, 2).Value = "Cronos"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(51, 4).Value = "0.0503"
$ws.Cells.Item(51, 5).Value = "  -1.58%  "
